# Actualización automática 2025-08-26 09:15:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("E19").Value = 74.36
$ws1.Range("E34").Value = "2 de 32"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F19").Value = 3462.21
$ws2.Range("F34").Value = 30837.34

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D4").Value = 203.94
$ws3.Range("E4").Value = 442.421575487259
$ws3.Range("F4").Value = 0.3155199933508734
$ws3.Range("D19").Value = 31355.55
$ws3.Range("E19").Value = 753.7310755578726
$ws3.Range("F19").Value = 0.9765260681550537
